# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G on Sheet1 is labeled "K" (strikeouts). The previously saved values
# were computed from a "Strike#" (total pitch/strike count) style stat; this
# regenerates that column using the correct K-based values for each game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$kValues = [ordered]@{
    2  = 1
    3  = 2
    4  = 4
    5  = 3
    6  = 2
    7  = 7
    8  = 2
    9  = 3
    10 = 6
    11 = 1
    12 = 1
    13 = 3
    14 = 3
    15 = 2
    16 = 3
    17 = 2
    18 = 1
    19 = 2
    20 = 1
    21 = 2
    22 = 4
    23 = 4
    24 = 1
    25 = 1
    26 = 3
    27 = 2
    28 = 2
    29 = 3
    30 = 1
    31 = 2
    32 = 1
    33 = 1
    34 = 0
    35 = 1
    36 = 3
    37 = 1
    38 = 2
    40 = 2
    41 = 2
    42 = 1
    43 = 3
    44 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
